$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.903.66"
$ws.Range("E2").Value = "  +0.52%  "
$ws.Range("D3").Value = "3.495.24"
$ws.Range("E3").Value = "  +0.27%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").Value = "'593.78"
$ws.Range("E5").Value = "  +0.61%  "
$ws.Range("D6").Value = "'171.92"
$ws.Range("E6").Value = "  +1.63%  "
$ws.Range("E7").Value = "  +0.03%  "
$ws.Range("D8").Value = "'0.591"
$ws.Range("E8").Value = "  -0.07%  "
$ws.Range("D9").Value = "'0.132"
$ws.Range("E9").Value = "  +4.07%  "
$ws.Range("E10").Value = "  -1.21%  "
$ws.Range("D11").Value = "'0.432"
$ws.Range("E11").Value = "  -1.26%  "
$ws.Range("D12").Value = "4.100.28"
$ws.Range("E12").Value = "  +0.28%  "
$ws.Range("E13").Value = "  -0.32%  "
$ws.Range("D14").Value = "'29.14"
$ws.Range("E14").Value = "  +3.68%  "
$ws.Range("D15").Value = "66.929.01"
$ws.Range("E15").Value = "  +0.51%  "
$ws.Range("D16").Value = "'0.0000179"
$ws.Range("E16").Value = "  +0.42%  "
$ws.Range("D17").Value = "3.483.41"
$ws.Range("E17").Value = "  -0.36%  "
$ws.Range("E18").Value = "  -0.40%  "
$ws.Range("D19").Value = "'14.06"
$ws.Range("E19").Value = "  +0.09%  "
$ws.Range("D20").Value = "'394.95"
$ws.Range("E20").Value = "  +1.06%  "
$ws.Range("E21").Value = "  +0.81%  "
$ws.Range("D22").Value = "'73.35"
$ws.Range("E22").Value = "  +0.44%  "
$ws.Range("E23").Value = "  +0.09%  "
$ws.Range("D24").Value = "'0.535"
$ws.Range("E24").Value = "  +0.08%  "
$ws.Range("E25").Value = "  +0.21%  "
$ws.Range("D26").Value = "'10.25"
$ws.Range("E26").Value = "  +0.64%  "
$ws.Range("E27").Value = "  -0.11%  "
$ws.Range("D28").Value = "'0.995"
$ws.Range("E28").Value = "  -0.42%  "
$ws.Range("D29").Value = "'6.17"
$ws.Range("E29").Value = "  -2.01%  "
$ws.Range("E30").Value = "  -2.10%  "
$ws.Range("D31").Value = "'2.06"
$ws.Range("E31").Value = "  +0.00%  "
$ws.Range("D32").Value = "'23.71"
$ws.Range("E32").Value = "  +0.83%  "
$ws.Range("D33").Value = "'7.35"
$ws.Range("E33").Value = "  -0.91%  "
$ws.Range("D35").Value = "'162.96"
$ws.Range("E35").Value = "  +0.72%  "
$ws.Range("D36").Value = "'0.877"
$ws.Range("E36").Value = "  -1.30%  "
$ws.Range("E37").Value = "  -0.60%  "
$ws.Range("D38").Value = "'6.95"
$ws.Range("E38").Value = "  +3.43%  "
$ws.Range("D39").Value = "'4.65"
$ws.Range("E39").Value = "  +0.08%  "
$ws.Range("E40").Value = "  -0.05%  "
$ws.Range("D41").Value = "'27.28"
$ws.Range("E41").Value = "  +2.17%  "
$ws.Range("D42").Value = "2.831.07"
$ws.Range("E42").Value = "  +2.32%  "
$ws.Range("D43").Value = "'26.26"
$ws.Range("E43").Value = "  -0.70%  "
$ws.Range("D44").Value = "'42.81"
$ws.Range("E44").Value = "  -0.65%  "
$ws.Range("E45").Value = "  +2.95%  "
$ws.Range("E46").Value = "  -2.81%  "
$ws.Range("D47").Value = "'335.90"
$ws.Range("E47").Value = "  -2.87%  "
$ws.Range("D48").Value = "'34.64"
$ws.Range("E48").Value = "  +2.37%  "
$ws.Range("E49").Value = "  -1.10%  "
$ws.Range("D50").Value = "'6.42"
$ws.Range("E50").Value = "  -1.84%  "
$ws.Range("D51").Value = "'0.840"
$ws.Range("E51").Value = "  -4.76%  "
